# Update the Dnajb11-Prtg NATMI LR-pairs sheet with refreshed TPM-derived
# statistics (commit: "update scripts wuth new tpm").
#
# Columns A-F (Sending cluster, Ligand symbol, Receptor symbol, Target
# cluster, Ligand-expressing cells, Ligand detection rate) are untouched.
# Columns G-T (expression values/specificities/weights) are recomputed for
# every data row (2-16) using the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.220785333333334
$ws.Range("H2").Value = 21.662356
$ws.Range("I2").Value = 0.1140291962005235
$ws.Range("J2").Value = 0.1193238000203875
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01675366666666667
$ws.Range("N2").Value = 0.050261
$ws.Range("O2").Value = 0.01455609640253496
$ws.Range("P2").Value = 0.01645081268867244
$ws.Range("Q2").Value = 0.1209746305462222
$ws.Range("R2").Value = 1.088771674916
$ws.Range("S2").Value = 0.001659819972598393
$ws.Range("T2").Value = 0.001962973483436003

$ws.Range("G3").Value = 7.220785333333334
$ws.Range("H3").Value = 21.662356
$ws.Range("I3").Value = 0.1140291962005235
$ws.Range("J3").Value = 0.1193238000203875
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7365303333333334
$ws.Range("N3").Value = 2.209591
$ws.Range("O3").Value = 0.6399200096729796
$ws.Range("P3").Value = 0.7232161648112139
$ws.Range("Q3").Value = 5.318327428488446
$ws.Range("R3").Value = 47.86494685639601
$ws.Range("S3").Value = 0.0729695643356411
$ws.Range("T3").Value = 0.0862969010214449

$ws.Range("G4").Value = 7.220785333333334
$ws.Range("H4").Value = 21.662356
$ws.Range("I4").Value = 0.1140291962005235
$ws.Range("J4").Value = 0.1193238000203875
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3976885
$ws.Range("N4").Value = 0.795377
$ws.Range("O4").Value = 0.3455238939244856
$ws.Range("P4").Value = 0.2603330225001138
$ws.Range("Q4").Value = 2.871623288035334
$ws.Range("R4").Value = 17.229739728212
$ws.Range("S4").Value = 0.03939981189228405
$ws.Range("T4").Value = 0.03106392551550662

$ws.Range("G5").Value = 14.26882533333333
$ws.Range("H5").Value = 42.806476
$ws.Range("I5").Value = 0.2253304326850228
$ws.Range("J5").Value = 0.2357929756948652
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.01675366666666667
$ws.Range("N5").Value = 0.050261
$ws.Range("O5").Value = 0.01455609640253496
$ws.Range("P5").Value = 0.01645081268867244
$ws.Range("Q5").Value = 0.2390551433595556
$ws.Range("R5").Value = 2.151496290236
$ws.Range("S5").Value = 0.003279931500588106
$ws.Range("T5").Value = 0.00387898607646092

$ws.Range("G6").Value = 14.26882533333333
$ws.Range("H6").Value = 42.806476
$ws.Range("I6").Value = 0.2253304326850228
$ws.Range("J6").Value = 0.2357929756948652
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7365303333333334
$ws.Range("N6").Value = 2.209591
$ws.Range("O6").Value = 0.6399200096729796
$ws.Range("P6").Value = 0.7232161648112139
$ws.Range("Q6").Value = 10.50942267903511
$ws.Range("R6").Value = 94.58480411131602
$ws.Range("S6").Value = 0.1441934526634165
$ws.Range("T6").Value = 0.1705292915714642

$ws.Range("G7").Value = 14.26882533333333
$ws.Range("H7").Value = 42.806476
$ws.Range("I7").Value = 0.2253304326850228
$ws.Range("J7").Value = 0.2357929756948652
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3976885
$ws.Range("N7").Value = 0.795377
$ws.Range("O7").Value = 0.3455238939244856
$ws.Range("P7").Value = 0.2603330225001138
$ws.Range("Q7").Value = 5.674547743575333
$ws.Range("R7").Value = 34.047286461452
$ws.Range("S7").Value = 0.07785704852101827
$ws.Range("T7").Value = 0.06138469804694011

$ws.Range("G8").Value = 18.28299766666667
$ws.Range("H8").Value = 54.848993
$ws.Range("I8").Value = 0.2887214384344039
$ws.Range("J8").Value = 0.3021273527243128
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01675366666666667
$ws.Range("N8").Value = 0.050261
$ws.Range("O8").Value = 0.01455609640253496
$ws.Range("P8").Value = 0.01645081268867244
$ws.Range("Q8").Value = 0.3063072485747778
$ws.Range("R8").Value = 2.756765237173
$ws.Range("S8").Value = 0.004202657091329745
$ws.Range("T8").Value = 0.004970240487792138

$ws.Range("G9").Value = 18.28299766666667
$ws.Range("H9").Value = 54.848993
$ws.Range("I9").Value = 0.2887214384344039
$ws.Range("J9").Value = 0.3021273527243128
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7365303333333334
$ws.Range("N9").Value = 2.209591
$ws.Range("O9").Value = 0.6399200096729796
$ws.Range("P9").Value = 0.7232161648112139
$ws.Range("Q9").Value = 13.46598236576256
$ws.Range("R9").Value = 121.193841291863
$ws.Range("S9").Value = 0.1847586256757403
$ws.Range("T9").Value = 0.2185033853218424

$ws.Range("G10").Value = 18.28299766666667
$ws.Range("H10").Value = 54.848993
$ws.Range("I10").Value = 0.2887214384344039
$ws.Range("J10").Value = 0.3021273527243128
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3976885
$ws.Range("N10").Value = 0.795377
$ws.Range("O10").Value = 0.3455238939244856
$ws.Range("P10").Value = 0.2603330225001138
$ws.Range("Q10").Value = 7.270937917560166
$ws.Range("R10").Value = 43.625627505361
$ws.Range("S10").Value = 0.09976015566733389
$ws.Range("T10").Value = 0.07865372691467833

$ws.Range("G11").Value = 8.429387500000001
$ws.Range("H11").Value = 16.858775
$ws.Range("I11").Value = 0.1331151996238646
$ws.Range("J11").Value = 0.09286400319008276
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.01675366666666667
$ws.Range("N11").Value = 0.050261
$ws.Range("O11").Value = 0.01455609640253496
$ws.Range("P11").Value = 0.01645081268867244
$ws.Range("Q11").Value = 0.1412231483791667
$ws.Range("R11").Value = 0.8473388902750001
$ws.Range("S11").Value = 0.001937637678367658
$ws.Range("T11").Value = 0.001527688322000331

$ws.Range("G12").Value = 8.429387500000001
$ws.Range("H12").Value = 16.858775
$ws.Range("I12").Value = 0.1331151996238646
$ws.Range("J12").Value = 0.09286400319008276
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.7365303333333334
$ws.Range("N12").Value = 2.209591
$ws.Range("O12").Value = 0.6399200096729796
$ws.Range("P12").Value = 0.7232161648112139
$ws.Range("Q12").Value = 6.208499585170834
$ws.Range("R12").Value = 37.25099751102501
$ws.Range("S12").Value = 0.08518307983092403
$ws.Range("T12").Value = 0.06716074823614798

$ws.Range("G13").Value = 8.429387500000001
$ws.Range("H13").Value = 16.858775
$ws.Range("I13").Value = 0.1331151996238646
$ws.Range("J13").Value = 0.09286400319008276
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3976885
$ws.Range("N13").Value = 0.795377
$ws.Range("O13").Value = 0.3455238939244856
$ws.Range("P13").Value = 0.2603330225001138
$ws.Range("Q13").Value = 3.35227047079375
$ws.Range("R13").Value = 13.409081883175
$ws.Range("S13").Value = 0.04599448211457292
$ws.Range("T13").Value = 0.02417556663193445

$ws.Range("G14").Value = 15.12200866666667
$ws.Range("H14").Value = 45.36602600000001
$ws.Range("I14").Value = 0.2388037330561851
$ws.Range("J14").Value = 0.2498918683703518
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.01675366666666667
$ws.Range("N14").Value = 0.050261
$ws.Range("O14").Value = 0.01455609640253496
$ws.Range("P14").Value = 0.01645081268867244
$ws.Range("Q14").Value = 0.2533490925317778
$ws.Range("R14").Value = 2.280141832786
$ws.Range("S14").Value = 0.003476050159651055
$ws.Range("T14").Value = 0.004110924318983045

$ws.Range("G15").Value = 15.12200866666667
$ws.Range("H15").Value = 45.36602600000001
$ws.Range("I15").Value = 0.2388037330561851
$ws.Range("J15").Value = 0.2498918683703518
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.7365303333333334
$ws.Range("N15").Value = 2.209591
$ws.Range("O15").Value = 0.6399200096729796
$ws.Range("P15").Value = 0.7232161648112139
$ws.Range("Q15").Value = 11.13781808392956
$ws.Range("R15").Value = 100.240362755366
$ws.Range("S15").Value = 0.1528152871672576
$ws.Range("T15").Value = 0.1807258386603145

$ws.Range("G16").Value = 15.12200866666667
$ws.Range("H16").Value = 45.36602600000001
$ws.Range("I16").Value = 0.2388037330561851
$ws.Range("J16").Value = 0.2498918683703518
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3976885
$ws.Range("N16").Value = 0.795377
$ws.Range("O16").Value = 0.3455238939244856
$ws.Range("P16").Value = 0.2603330225001138
$ws.Range("Q16").Value = 6.013848943633668
$ws.Range("R16").Value = 36.083093661802
$ws.Range("S16").Value = 0.08251239572927649
$ws.Range("T16").Value = 0.06505510539105427
